# Updated 2D training schedules, no break screen
# Adds a new "break_on_off" column (L) to the schedule sheet, flagging the
# trials on which a break screen is shown (1) vs not (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column L
$ws.Cells.Item(1, 12).Value = "break_on_off"

# Trials (rows) on which the break screen is turned on
$breakRows = @(19, 37, 54)

# Fill column L for every data row (rows 2-73, trials 1-72)
for ($r = 2; $r -le 73; $r++) {
    if ($breakRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Reflect the user's selection of the newly added column in the saved view
[void]$ws.Range("L1:L73").Select()
